$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 84692300
$ws.Range("E8").Value = 82826900
$ws.Range("F8").Value = 90710100
$ws.Range("G8").Value = 88365400
$ws.Range("H8").Value = 87384700
$ws.Range("I8").Value = 81731300
$ws.Range("J8").Value = 87379600

$ws.Range("D9").Value = 62073400
$ws.Range("E9").Value = 61315400
$ws.Range("F9").Value = 67430000
$ws.Range("G9").Value = 65072000
$ws.Range("H9").Value = 64850700
$ws.Range("I9").Value = 60965700
$ws.Range("J9").Value = 65801900

$ws.Range("D10").Value = 22618900
$ws.Range("E10").Value = 21511500
$ws.Range("F10").Value = 23280100
$ws.Range("G10").Value = 23293300
$ws.Range("H10").Value = 22534000
$ws.Range("I10").Value = 20765600
$ws.Range("J10").Value = 21577700

$ws.Range("H12").Value = 3176900
$ws.Range("I12").Value = 3085400
$ws.Range("J12").Value = 3729100

$ws.Range("H14").Value = 561200
$ws.Range("I14").Value = 485500
$ws.Range("J14").Value = 496600

$ws.Range("D17").Value = 78232000
$ws.Range("E17").Value = 77517600
$ws.Range("F17").Value = 84970900
$ws.Range("G17").Value = 82567800
$ws.Range("H17").Value = 81917300
$ws.Range("I17").Value = 78401700
$ws.Range("J17").Value = 84149200

$ws.Range("D18").Value = 6460300
$ws.Range("E18").Value = 5309300
$ws.Range("F18").Value = 5739200
$ws.Range("G18").Value = 5797600
$ws.Range("H18").Value = 5467400
$ws.Range("I18").Value = 3329600
$ws.Range("J18").Value = 3230400

$ws.Range("D20").Value = -501200
$ws.Range("E20").Value = -896800
$ws.Range("F20").Value = -830200
$ws.Range("G20").Value = -856400
$ws.Range("H20").Value = 1145500
$ws.Range("I20").Value = 26400
$ws.Range("J20").Value = 2065900

$ws.Range("D21").Value = 9257100
$ws.Range("E21").Value = 8169900
$ws.Range("F21").Value = 9504500
$ws.Range("G21").Value = 9294400
$ws.Range("H21").Value = 11870000
$ws.Range("I21").Value = 7139100
$ws.Range("J21").Value = 9637300

$ws.Range("D22").Value = 185700
$ws.Range("E22").Value = 171900
$ws.Range("F22").Value = 235000
$ws.Range("G22").Value = 249400
$ws.Range("H22").Value = 479300
$ws.Range("I22").Value = 241400
$ws.Range("J22").Value = 254400

$ws.Range("D23").Value = 5773400
$ws.Range("E23").Value = 4240600
$ws.Range("F23").Value = 4674000
$ws.Range("G23").Value = 4691700
$ws.Range("H23").Value = 6133600
$ws.Range("I23").Value = 3114600
$ws.Range("J23").Value = 5041900

$ws.Range("D24").Value = 1190600
$ws.Range("E24").Value = 1131000
$ws.Range("F24").Value = 1493500
$ws.Range("G24").Value = 1103600
$ws.Range("H24").Value = 1324700
$ws.Range("I24").Value = 965600
$ws.Range("J24").Value = 1310100

$ws.Range("D26").Value = 4582700
$ws.Range("E26").Value = 3109600
$ws.Range("F26").Value = 3180600
$ws.Range("G26").Value = 3588100
$ws.Range("H26").Value = 4808900
$ws.Range("I26").Value = 2149000
$ws.Range("J26").Value = 3731800

$ws.Range("D27").Value = 3426200
$ws.Range("E27").Value = 2144400
$ws.Range("F27").Value = 2072300
$ws.Range("G27").Value = 2449700
$ws.Range("H27").Value = 2908900
$ws.Range("I27").Value = 1584900
$ws.Range("J27").Value = 3138500

$ws.Range("D29").Value = -144800
$ws.Range("E29").Value = -53800
$ws.Range("F29").Value = -516000
$ws.Range("G29").Value = -483600
$ws.Range("H29").Value = -62900

$ws.Range("D32").Value = 501200
$ws.Range("E32").Value = 896800
$ws.Range("F32").Value = 830200
$ws.Range("G32").Value = 856400
$ws.Range("H32").Value = -1145500
$ws.Range("I32").Value = -26400
$ws.Range("J32").Value = -2065900

$ws.Range("D33").Value = 3281400
$ws.Range("E33").Value = 2090600
$ws.Range("F33").Value = 1556300
$ws.Range("G33").Value = 1966000
$ws.Range("H33").Value = 2846000
$ws.Range("I33").Value = 1584900
$ws.Range("J33").Value = 3138500

$ws.Range("D35").Value = 3281400
$ws.Range("E35").Value = 2090600
$ws.Range("F35").Value = 1556300
$ws.Range("G35").Value = 1966000
$ws.Range("H35").Value = 2846000
$ws.Range("I35").Value = 1584900
$ws.Range("J35").Value = 3138500

$ws.Range("D41").Value = 6309600
$ws.Range("E41").Value = 6917800
$ws.Range("F41").Value = 6321800
$ws.Range("G41").Value = 6343400
$ws.Range("H41").Value = 10114600
$ws.Range("I41").Value = 4769800
$ws.Range("J41").Value = 5601000

$ws.Range("D42").Value = 3374800
$ws.Range("H42").Value = 554900
$ws.Range("I42").Value = 305600
$ws.Range("J42").Value = 882600

$ws.Range("D43").Value = 22612800
$ws.Range("E43").Value = 22378600
$ws.Range("F43").Value = 30117000
$ws.Range("G43").Value = 28994900
$ws.Range("H43").Value = 51462500
$ws.Range("I43").Value = 21892900
$ws.Range("J43").Value = 21185000

$ws.Range("D44").Value = 12432100
$ws.Range("E44").Value = 11082200
$ws.Range("F44").Value = 11750700
$ws.Range("G44").Value = 13181400
$ws.Range("H44").Value = 24824300
$ws.Range("I44").Value = 12994100
$ws.Range("J44").Value = 12775800

$ws.Range("D45").Value = 1842900
$ws.Range("E45").Value = 4844900
$ws.Range("F45").Value = 4898400
$ws.Range("G45").Value = 4657400
$ws.Range("H45").Value = 12175800
$ws.Range("I45").Value = 6956500
$ws.Range("J45").Value = 6221800

$ws.Range("D46").Value = 46572300
$ws.Range("E46").Value = 45223600
$ws.Range("F46").Value = 53087900
$ws.Range("G46").Value = 53177000
$ws.Range("H46").Value = 47569200
$ws.Range("I46").Value = 46918800
$ws.Range("J46").Value = 46666200

$ws.Range("D47").Value = 13196900
$ws.Range("E47").Value = 12755000
$ws.Range("F47").Value = 18142700
$ws.Range("G47").Value = 19267500
$ws.Range("H47").Value = 29805900
$ws.Range("I47").Value = 8256800
$ws.Range("J47").Value = 8587100

$ws.Range("D48").Value = 19208400
$ws.Range("E48").Value = 18065600
$ws.Range("F48").Value = 22602000
$ws.Range("G48").Value = 22351400
$ws.Range("H48").Value = 41593300
$ws.Range("I48").Value = 20610900
$ws.Range("J48").Value = 18310900

$ws.Range("D49").Value = 9531500
$ws.Range("E49").Value = 8309600
$ws.Range("F49").Value = 9676400
$ws.Range("G49").Value = 8439600
$ws.Range("H49").Value = 13503200
$ws.Range("I49").Value = 6376800
$ws.Range("J49").Value = 5514100

$ws.Range("D52").Value = 2854500
$ws.Range("E52").Value = 3008000
$ws.Range("F52").Value = 9952000
$ws.Range("G52").Value = 9165500
$ws.Range("H52").Value = 15885900
$ws.Range("I52").Value = 6512100
$ws.Range("J52").Value = 6065300

$ws.Range("D54").Value = 91363700
$ws.Range("E54").Value = 87361800
$ws.Range("F54").Value = 113461000
$ws.Range("G54").Value = 112401000
$ws.Range("H54").Value = 100328000
$ws.Range("I54").Value = 88675400
$ws.Range("J54").Value = 85143500

$ws.Range("D57").Value = 13894300
$ws.Range("E57").Value = 12676200
$ws.Range("F57").Value = 13125300
$ws.Range("G57").Value = 12895800
$ws.Range("H57").Value = 24384500
$ws.Range("I57").Value = 11163200
$ws.Range("J57").Value = 11767900

$ws.Range("D58").Value = 2157200
$ws.Range("E58").Value = 3494800
$ws.Range("F58").Value = 13767300
$ws.Range("G58").Value = 13209400
$ws.Range("H58").Value = 22404400
$ws.Range("I58").Value = 8682300
$ws.Range("J58").Value = 9370700

$ws.Range("D59").Value = 18258800
$ws.Range("E59").Value = 17465600
$ws.Range("F59").Value = 18255000
$ws.Range("G59").Value = 17101300
$ws.Range("H59").Value = 30875600
$ws.Range("I59").Value = 15990300
$ws.Range("J59").Value = 16023700

$ws.Range("D60").Value = 34310400
$ws.Range("E60").Value = 33636600
$ws.Range("F60").Value = 45147700
$ws.Range("G60").Value = 43206500
$ws.Range("H60").Value = 38876000
$ws.Range("I60").Value = 35835800
$ws.Range("J60").Value = 37162300

$ws.Range("D61").Value = 7337400
$ws.Range("E61").Value = 7141700
$ws.Range("F61").Value = 19857900
$ws.Range("G61").Value = 20011600
$ws.Range("H61").Value = 16423500
$ws.Range("I61").Value = 12743200
$ws.Range("J61").Value = 12510400

$ws.Range("D62").Value = 8930400
$ws.Range("E62").Value = 9546700
$ws.Range("F62").Value = 11160300
$ws.Range("G62").Value = 10343900
$ws.Range("H62").Value = 21010500
$ws.Range("I62").Value = 11355700
$ws.Range("J62").Value = 10393900

$ws.Range("D66").Value = 61730400
$ws.Range("E66").Value = 60539400
$ws.Range("F66").Value = 88736000
$ws.Range("G66").Value = 85802700
$ws.Range("H66").Value = 76203000
$ws.Range("I66").Value = 69849100
$ws.Range("J66").Value = 69126600

$ws.Range("D72").Value = 19032800
$ws.Range("E72").Value = 16213900
$ws.Range("F72").Value = 14552200
$ws.Range("G72").Value = 13356800
$ws.Range("H72").Value = 25902900
$ws.Range("I72").Value = 12391300
$ws.Range("J72").Value = 11228700

$ws.Range("D76").Value = 29633300
$ws.Range("E76").Value = 26822400
$ws.Range("F76").Value = 24725100
$ws.Range("G76").Value = 26598200
$ws.Range("H76").Value = 24124700
$ws.Range("I76").Value = 18826300
$ws.Range("J76").Value = 16016900

$ws.Range("D81").Value = 3281400
$ws.Range("E81").Value = 2090600
$ws.Range("F81").Value = 1556300
$ws.Range("G81").Value = 1966000
$ws.Range("H81").Value = 2846000
$ws.Range("I81").Value = 1584900
$ws.Range("J81").Value = 3138500

$ws.Range("D83").Value = 3294500
$ws.Range("E83").Value = 3753300
$ws.Range("F83").Value = 4590400
$ws.Range("G83").Value = 4348400
$ws.Range("H83").Value = 5251300
$ws.Range("I83").Value = 3778900
$ws.Range("J83").Value = 4336200

$ws.Range("D89").Value = 6573600
$ws.Range("E89").Value = 5691400
$ws.Range("F89").Value = 7342500
$ws.Range("G89").Value = 4084500
$ws.Range("H89").Value = 2773300
$ws.Range("I89").Value = 5274900
$ws.Range("J89").Value = 4042300

$ws.Range("D91").Value = -3188000
$ws.Range("E91").Value = -5505900
$ws.Range("F91").Value = -8216600
$ws.Range("G91").Value = -7253400
$ws.Range("H91").Value = -14879600
$ws.Range("I91").Value = -6485000
$ws.Range("J91").Value = -5485900

$ws.Range("D94").Value = -4287900
$ws.Range("E94").Value = -3055100
$ws.Range("F94").Value = -6606400
$ws.Range("G94").Value = -5537400
$ws.Range("H94").Value = -4973600
$ws.Range("I94").Value = -5003300
$ws.Range("J94").Value = -1768100

$ws.Range("D96").Value = -610800
$ws.Range("E96").Value = -523700
$ws.Range("F96").Value = -523500
$ws.Range("G96").Value = -501200
$ws.Range("H96").Value = -435700
$ws.Range("I96").Value = -420900
$ws.Range("J96").Value = -244200

$ws.Range("D100").Value = -2905900
$ws.Range("E100").Value = -1894200
$ws.Range("F100").Value = -239300
$ws.Range("G100").Value = 2108200
$ws.Range("H100").Value = 2068700
$ws.Range("I100").Value = -1631200
$ws.Range("J100").Value = -1517300

$ws.Range("D101").Value = 12100
$ws.Range("E101").Value = -146100
$ws.Range("F101").Value = -518400
$ws.Range("G101").Value = 619800
$ws.Range("H101").Value = 468800
$ws.Range("I101").Value = 528400
$ws.Range("J101").Value = -171500

$ws.Range("D102").Value = -608200
$ws.Range("E102").Value = 596000
$ws.Range("F102").Value = -21600
$ws.Range("G102").Value = 1275100
$ws.Range("H102").Value = 337200
$ws.Range("I102").Value = -831200
$ws.Range("J102").Value = 585500
